$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: update live aggregate metrics now that trade #53 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.49              # Current Capital
$summary.Range("B4").Value = -2.5                 # Total P&L $
$summary.Range("B5").Value = -0.9399999999999999  # Total P&L %
$summary.Range("B6").Value = 53                   # Total Trades
$summary.Range("B7").Value = 21                   # Winning Trades
$summary.Range("B9").Value = 39.62                # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: refresh the MarketMaking strategy row (row 4).
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.48999999999999  # Capital
$status.Range("D4").Value = 53                 # Trades
$status.Range("E4").Value = -2.5               # P&L $
$status.Range("F4").Value = -2.51              # P&L %
$status.Range("G4").Value = 39.62              # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade (#53) as row 54 to both the "All Trades"
# and "MarketMaking" trade logs.
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(54, 1).Value = 53                 # A: Trade #

    # Force the Date column to stay plain text (matches the rest of the
    # column) instead of Excel auto-converting "2026-02-17" into a date.
    $ws.Cells.Item(54, 2).NumberFormat = "@"
    $ws.Cells.Item(54, 2).Value = "2026-02-17"        # B: Date

    $ws.Cells.Item(54, 3).Value = "13:29:11"          # C: Time
    $ws.Cells.Item(54, 4).Value = "MarketMaking"      # D: Strategy
    $ws.Cells.Item(54, 5).Value = "UP"                # E: Side
    $ws.Cells.Item(54, 6).Value = 0.91                # F: Entry Price
    $ws.Cells.Item(54, 7).Value = 0.93                # G: Exit Price
    $ws.Cells.Item(54, 8).Value = "CLOSED"            # H: Status
    $ws.Cells.Item(54, 9).Value = 2.1978              # I: P&L %
    $ws.Cells.Item(54, 10).Value = 0.02               # J: P&L $
    $ws.Cells.Item(54, 11).Value = 97.48999999999999  # K: Capital After
    $ws.Cells.Item(54, 12).Value = 0                  # L: Entry Slippage (bps)
    $ws.Cells.Item(54, 13).Value = 0                  # M: Exit Slippage (bps)
    $ws.Cells.Item(54, 14).Value = 0.6                # N: Confidence
    $ws.Cells.Item(54, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
    $ws.Cells.Item(54, 16).Value = "early_exit"       # P: Exit Reason
    $ws.Cells.Item(54, 17).Value = 0.1                # Q: Duration (min)
}
